# "icons on each button vers2"
# Adds a new "icons" worksheet after "comments", updates the remembered
# selection on each existing sheet, tweaks one cell style on the URL sheet,
# and leaves the new "icons" sheet as the active tab/selection.

$wb = $excel.ActiveWorkbook

# --- fieldnames: move remembered selection from A2 to B2 -------------------
$wsFieldnames = $wb.Worksheets.Item("fieldnames")
$wsFieldnames.Activate()
$wsFieldnames.Range("B2").Select()

# --- URL: move remembered selection from C20 to B2, restyle C6 -------------
$wsUrl = $wb.Worksheets.Item("URL")
$wsUrl.Activate()
$wsUrl.Range("B6").Copy()
$wsUrl.Range("C6").PasteSpecial(-4122)
$wsUrl.Range("B2").Select()

# --- color: move remembered selection from B6 to B3 -------------------------
$wsColor = $wb.Worksheets.Item("color")
$wsColor.Activate()
$wsColor.Range("B3").Select()

# --- comments: move remembered selection from B2 to B3 ----------------------
$wsComments = $wb.Worksheets.Item("comments")
$wsComments.Activate()
$wsComments.Range("B3").Select()

# --- icons: brand new sheet added after "comments" --------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsIcons = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsIcons.Name = "icons"

$wsIcons.Range("A2").Value = "globe.png"
$wsIcons.Range("B2").Value = "BCN_Logo3.png"
$wsIcons.Range("C2").Value = "globe.png"
$wsIcons.Range("D2").Value = "info-circle.png"

$wsIcons.Range("A3").Value = "globe.png"
$wsIcons.Range("B3").Value = "BCN_Logo3.png"
$wsIcons.Range("C3").Value = "globe.png"

$wsIcons.Range("B4").Value = "info-circle.png"
$wsIcons.Range("C4").Value = "info-circle_red.png"

$wsIcons.Range("B5").Value = "info-circle.png"
$wsIcons.Range("C5").Value = "info-circle.png"

$wsIcons.Range("B6").Value = "info-circle.png"
$wsIcons.Range("C6").Value = "info-circle.png"

$wsIcons.Range("B7").Value = "info-circle.png"
$wsIcons.Range("C7").Value = "info-circle.png"

$wsIcons.Range("C8").Value = "info-circle.png"
$wsIcons.Range("C9").Value = "info-circle.png"

$wsIcons.Activate()
$wsIcons.Range("C5").Select()
